# Helper: write a value into a cell while forcing it to stay a *text* cell
# (the source workbook stores every "Calculated / Unit attributes" value as
# a string, even when it looks numeric/boolean, e.g. "0.0", "False"). Excel's
# normal Value setter auto-coerces such literals to Number/Boolean, so we use
# a leading quote to force text entry, then clear the resulting "quote
# prefix" cell format so no stray style is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) X101_TBP_Extraction!B19 - reorder species list (nonextract_to_aq)
# ---------------------------------------------------------------------------
$wsExtract = $wb.Worksheets.Item("X101_TBP_Extraction")
$wsExtract.Range("B19").Value = "HNO3, Gd(NO3)3, H2O, Eu(NO3)3, Nd(NO3)3, I_aq, Sr(NO3)2, CsNO3, Sm(NO3)3"

# ---------------------------------------------------------------------------
# 2) X102_AHA_Strip!B18 - reorder species list (nontransfer_keep_in_aq)
# ---------------------------------------------------------------------------
$wsStrip = $wb.Worksheets.Item("X102_AHA_Strip")
$wsStrip.Range("B18").Value = "HNO3, H2O, AHA"

# ---------------------------------------------------------------------------
# 3) TSA101A_ColA - insert "captured_cycle_mol" + "regen_source" attribute
#    rows, update I_desorb_mol_s and print_diagnostics values.
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("TSA101A_ColA")

# I_desorb_mol_s (row 11) now carries the desorbed amount.
Set-TextValue $wsA.Range("B11") "0.00108381924"

# Insert new row 18 "captured_cycle_mol" (pushes iodine_name.. down by one).
$wsA.Rows(18).Insert()
$wsA.Range("A18").Value = "captured_cycle_mol"
Set-TextValue $wsA.Range("B18") "15.606997056"

# print_diagnostics (now row 22) flips False -> True.
Set-TextValue $wsA.Range("B22") "True"

# Insert new row 23 "regen_source" (pushes regen_yI2_max.. down by one). Left
# blank, like the other section-header spacer cells on this sheet.
$wsA.Rows(23).Insert()
$wsA.Range("A23").Value = "regen_source"

# ---------------------------------------------------------------------------
# 4) TSA101B_ColB - same row insertions, plus regen_source value and the
#    updated regen_in / regen_out flow numbers.
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("TSA101B_ColB")

# Insert new row 18 "captured_cycle_mol" (pushes iodine_name.. down by one).
$wsB.Rows(18).Insert()
$wsB.Range("A18").Value = "captured_cycle_mol"
Set-TextValue $wsB.Range("B18") "0.0"

# print_diagnostics (now row 22) flips False -> True.
Set-TextValue $wsB.Range("B22") "True"

# Insert new row 23 "regen_source" (pushes regen_yI2_max.. down by one) and
# populate it with the repr() of the partner column object.
$wsB.Rows(23).Insert()
$wsB.Range("A23").Value = "regen_source"
Set-TextValue $wsB.Range("B23") "<process_sim.unitops.TSA.IdealTSAColumnEMM17 object at 0x0000018955AC5160>"

# regen_in (now row 31) F_total/m_dot split in half between the two columns.
$wsB.Range("F31").Value = 0.009754373160000001
$wsB.Range("G31").Value = 0.2814199084648225
